# Scheduled-runner refresh: updates currentAveragePrice / Leve price / profit
# columns (H,I,J,K,L,M,N) with freshly pulled market data, per item row, across
# all eight crafting-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Cells.Item(2, 8).Value = 1407.9166
$ws.Cells.Item(2, 9).Value = 1494
$ws.Cells.Item(2, 10).Value = 1149.6666
$ws.Cells.Item(2, 11).Value = 1494
$ws.Cells.Item(2, 12).Value = 1149.6666
$ws.Cells.Item(2, 13).Value = -1381
$ws.Cells.Item(2, 14).Value = -1375.6666

# Row 33
$ws.Cells.Item(33, 8).Value = 4071.75
$ws.Cells.Item(33, 9).Value = 4926.6875
$ws.Cells.Item(33, 10).Value = 652
$ws.Cells.Item(33, 11).Value = 4926.6875
$ws.Cells.Item(33, 12).Value = 652
$ws.Cells.Item(33, 13).Value = -4697.6875
$ws.Cells.Item(33, 14).Value = -1110

# Row 116
$ws.Cells.Item(116, 8).Value = 5001.778
$ws.Cells.Item(116, 9).Value = 3766.1428
$ws.Cells.Item(116, 10).Value = 5788.091
$ws.Cells.Item(116, 11).Value = 3766.1428
$ws.Cells.Item(116, 12).Value = 5788.091
$ws.Cells.Item(116, 13).Value = -324.1428000000001
$ws.Cells.Item(116, 14).Value = -12672.091

# Row 129
$ws.Cells.Item(129, 8).Value = 62502988
$ws.Cells.Item(129, 9).Value = 166667780
$ws.Cells.Item(129, 11).Value = 500003340
$ws.Cells.Item(129, 13).Value = -499998340

# Row 131
$ws.Cells.Item(131, 8).Value = 4047.9167
$ws.Cells.Item(131, 9).Value = 2479.4119
$ws.Cells.Item(131, 11).Value = 7438.2357
$ws.Cells.Item(131, 13).Value = -2398.2357

# Row 135
$ws.Cells.Item(135, 8).Value = 913.2273
$ws.Cells.Item(135, 9).Value = 890.35297
$ws.Cells.Item(135, 11).Value = 8013.17673
$ws.Cells.Item(135, 13).Value = -5478.17673

# Row 137
$ws.Cells.Item(137, 8).Value = 114072.19
$ws.Cells.Item(137, 9).Value = 297851.5
$ws.Cells.Item(137, 10).Value = 3804.6
$ws.Cells.Item(137, 11).Value = 893554.5
$ws.Cells.Item(137, 12).Value = 11413.8
$ws.Cells.Item(137, 13).Value = -891004.5
$ws.Cells.Item(137, 14).Value = -16513.8

# Row 138
$ws.Cells.Item(138, 8).Value = 4245.9062
$ws.Cells.Item(138, 9).Value = 2916.5
$ws.Cells.Item(138, 10).Value = 4552.6924
$ws.Cells.Item(138, 11).Value = 8749.5
$ws.Cells.Item(138, 12).Value = 13658.0772
$ws.Cells.Item(138, 13).Value = -3609.5
$ws.Cells.Item(138, 14).Value = -23938.0772

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 4274957
$ws.Cells.Item(2, 10).Value = 2599.3333
$ws.Cells.Item(2, 12).Value = 2599.3333
$ws.Cells.Item(2, 14).Value = -2825.3333

# Row 45
$ws.Cells.Item(45, 10).Value = 8803.5
$ws.Cells.Item(45, 12).Value = 8803.5
$ws.Cells.Item(45, 14).Value = -9557.5

# Row 64
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).Value = $null  # L64
$ws.Cells.Item(64, 14).Value = 0

# Row 67
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).Value = $null  # L67
$ws.Cells.Item(67, 14).Value = 0

# Row 116
$ws.Cells.Item(116, 8).Value = 4274957
$ws.Cells.Item(116, 10).Value = 2599.3333
$ws.Cells.Item(116, 12).Value = 2599.3333
$ws.Cells.Item(116, 14).Value = -7187.3333

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 4274957
$ws.Cells.Item(3, 10).Value = 2599.3333
$ws.Cells.Item(3, 12).Value = 2599.3333
$ws.Cells.Item(3, 14).Value = -2827.3333

# Row 22
$ws.Cells.Item(22, 8).Value = 1212.6
$ws.Cells.Item(22, 9).Value = 1140.6875
$ws.Cells.Item(22, 10).Value = 1500.25
$ws.Cells.Item(22, 11).Value = 1140.6875
$ws.Cells.Item(22, 12).Value = 1500.25
$ws.Cells.Item(22, 13).Value = -967.6875
$ws.Cells.Item(22, 14).Value = -1846.25

# Row 62
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 12).Value = $null  # L62
$ws.Cells.Item(62, 14).Value = 0

# Row 65
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 12).Value = $null  # L65
$ws.Cells.Item(65, 14).Value = 0

# Row 107
$ws.Cells.Item(107, 8).Value = 2464870.8
$ws.Cells.Item(107, 10).Value = 3889.9
$ws.Cells.Item(107, 12).Value = 3889.9
$ws.Cells.Item(107, 14).Value = -7729.9

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Cells.Item(22, 8).Value = 570
$ws.Cells.Item(22, 9).Value = 539.1429000000001
$ws.Cells.Item(22, 10).Value = 1002
$ws.Cells.Item(22, 11).Value = 539.1429000000001
$ws.Cells.Item(22, 12).Value = 1002
$ws.Cells.Item(22, 13).Value = -189.1429000000001
$ws.Cells.Item(22, 14).Value = -1702

# Row 31
$ws.Cells.Item(31, 8).Value = 15090.021
$ws.Cells.Item(31, 9).Value = 8430.200000000001
$ws.Cells.Item(31, 10).Value = 16387.389
$ws.Cells.Item(31, 11).Value = 8430.200000000001
$ws.Cells.Item(31, 12).Value = 16387.389
$ws.Cells.Item(31, 13).Value = -8135.200000000001
$ws.Cells.Item(31, 14).Value = -16977.389

# Row 34
$ws.Cells.Item(34, 8).Value = 15090.021
$ws.Cells.Item(34, 9).Value = 8430.200000000001
$ws.Cells.Item(34, 10).Value = 16387.389
$ws.Cells.Item(34, 11).Value = 8430.200000000001
$ws.Cells.Item(34, 12).Value = 16387.389
$ws.Cells.Item(34, 13).Value = -8228.200000000001
$ws.Cells.Item(34, 14).Value = -16791.389

# Row 64
$ws.Cells.Item(64, 8).Value = 59995
$ws.Cells.Item(64, 10).Value = 59995
$ws.Cells.Item(64, 12).Value = 59995
$ws.Cells.Item(64, 14).Value = -60491

# Row 67
$ws.Cells.Item(67, 8).Value = 59995
$ws.Cells.Item(67, 10).Value = 59995
$ws.Cells.Item(67, 12).Value = 59995
$ws.Cells.Item(67, 14).Value = -61711

# Row 99
$ws.Cells.Item(99, 8).Value = 6145.8
$ws.Cells.Item(99, 9).Value = 5122.25
$ws.Cells.Item(99, 11).Value = 5122.25
$ws.Cells.Item(99, 13).Value = -3624.25

# Row 126
$ws.Cells.Item(126, 8).Value = 6145.8
$ws.Cells.Item(126, 9).Value = 5122.25
$ws.Cells.Item(126, 11).Value = 15366.75
$ws.Cells.Item(126, 13).Value = -12896.75

# Row 134
$ws.Cells.Item(134, 8).Value = 8493.714
$ws.Cells.Item(134, 9).Value = 6193.1665
$ws.Cells.Item(134, 11).Value = 18579.4995
$ws.Cells.Item(134, 13).Value = -16044.4995

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Cells.Item(21, 8).Value = 6950
$ws.Cells.Item(21, 10).Value = 9750
$ws.Cells.Item(21, 12).Value = 9750
$ws.Cells.Item(21, 14).Value = -10096

# Row 30
$ws.Cells.Item(30, 8).Value = 6950
$ws.Cells.Item(30, 10).Value = 9750
$ws.Cells.Item(30, 12).Value = 9750
$ws.Cells.Item(30, 14).Value = -9960

# Row 47
$ws.Cells.Item(47, 8).Value = 49997.5
$ws.Cells.Item(47, 10).Value = 49997.5
$ws.Cells.Item(47, 12).Value = 49997.5
$ws.Cells.Item(47, 14).Value = -51133.5

$ws = $wb.Worksheets.Item("LTW")
# Row 23
$ws.Cells.Item(23, 8).Value = 5766.3335
$ws.Cells.Item(23, 9).Value = 5766.3335
$ws.Cells.Item(23, 11).Value = 5766.3335
$ws.Cells.Item(23, 13).Value = -5536.3335

# Row 97
$ws.Cells.Item(97, 8).Value = 18371.5
$ws.Cells.Item(97, 10).Value = 18371.5
$ws.Cells.Item(97, 12).Value = 18371.5
$ws.Cells.Item(97, 14).Value = -20353.5

$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Cells.Item(64, 8).Value = 49500
$ws.Cells.Item(64, 10).Value = 49500
$ws.Cells.Item(64, 12).Value = 49500
$ws.Cells.Item(64, 14).Value = -49996

# Row 67
$ws.Cells.Item(67, 8).Value = 49500
$ws.Cells.Item(67, 10).Value = 49500
$ws.Cells.Item(67, 12).Value = 49500
$ws.Cells.Item(67, 14).Value = -51216

# Row 127
$ws.Cells.Item(127, 8).Value = 40000
$ws.Cells.Item(127, 10).Value = 40000
$ws.Cells.Item(127, 12).Value = 40000
$ws.Cells.Item(127, 14).Value = -49920

